$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Color helpers (Excel COM colors are decimal BGR: R + G*256 + B*65536)
#   Font color FF67595E  -> 6183271
#   Fill color FFF6EEE1  -> 14806774
# ---------------------------------------------------------------------------
$fontColor = 6183271
$fillColor = 14806774

# ---------------------------------------------------------------------------
# Row 1: recipe title "Easy Chicken Pot Pie" -- build the "title" style here.
# (Color is set before Bold/Size so the "color-only" and "bold+color"
# intermediate font states are shared with the sub-header/body styles built
# below instead of each minting their own throwaway font entries.)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Easy Chicken Pot Pie"
$ws.Range("A1").Font.Color = $fontColor
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 18
$ws.Range("A1").Interior.Color = $fillColor
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108
$ws.Range("A1").WrapText = $true

# Row 2: blank spacer row -- just the fill color.
$ws.Range("A2").Interior.Color = $fillColor

# Row 3: "Ingredients" sub-header -- build the "sub-header" style here.
$ws.Range("A3").Value = "Ingredients"
$ws.Range("A3").Font.Color = $fontColor
$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").Font.Size = 14
$ws.Range("A3").Interior.Color = $fillColor
$ws.Range("A3").HorizontalAlignment = -4131
$ws.Range("A3").VerticalAlignment = -4108

# Row 4: ingredients list -- build the "body" style here.
$ws.Range("A4").Value = "1.0 can boned chicken, 1.0 can vegetables, 1.0 can cream of chicken, 1.0 cup bisquick, 1.0 cup water"
$ws.Range("A4").Font.Color = $fontColor
$ws.Range("A4").Font.Size = 12
$ws.Range("A4").Interior.Color = $fillColor
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("A4").VerticalAlignment = -4108
$ws.Range("A4").WrapText = $true

# Row 5: blank spacer row -- reuse the spacer style via copy/paste so no new
# styles get minted.
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)

# Row 6: "Instructions" sub-header -- reuse the sub-header style.
$ws.Range("A6").Value = "Instructions"
$ws.Range("A3").Copy()
$ws.Range("A6").PasteSpecial(-4122)

# Row 7: instructions text (multi-line) -- reuse the body style.
$ws.Range("A7").Value = "Step 1: Mix first 3 ingredients and pour into a greased baking dish.
Step 2: Mix Bisquick with water, pour over first mixture and bake in a 400 degree oven for 45 to 50 minutes. (You may use left over chicken or turkey.)"
$ws.Range("A4").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Rows(7).AutoFit()

# ---------------------------------------------------------------------------
# Row 9: second recipe title -- reuse the title style.
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "Skillet Roasted Chicken & Potatoes"
$ws.Range("A1").Copy()
$ws.Range("A9").PasteSpecial(-4122)

# Row 10: blank spacer row.
$ws.Range("A2").Copy()
$ws.Range("A10").PasteSpecial(-4122)

# Row 11: "Ingredients" sub-header.
$ws.Range("A11").Value = "Ingredients"
$ws.Range("A3").Copy()
$ws.Range("A11").PasteSpecial(-4122)

# Row 12: ingredients list.
$ws.Range("A12").Value = "3.0 Tbsps olive oil, 2.0 tsps thyme, 1.5 tsps paprika, 4.0 servings salt & pepper, 4.0 lb chicken, 2.0 lb yukon gold potatoes"
$ws.Range("A4").Copy()
$ws.Range("A12").PasteSpecial(-4122)

# Row 13: blank spacer row.
$ws.Range("A2").Copy()
$ws.Range("A13").PasteSpecial(-4122)

# Row 14: "Instructions" sub-header.
$ws.Range("A14").Value = "Instructions"
$ws.Range("A3").Copy()
$ws.Range("A14").PasteSpecial(-4122)

# Row 15: instructions text (multi-line).
$ws.Range("A15").Value = "Step 1: Adjust oven rack to lower-middle position and preheat oven to 400 degrees f. In a small bowl combine 2 T. olive oil, thyme, smoked paprika, 1 t. salt and 1/2 t. pepper. Carefully run your fingers under the skin covering the chicken breast, separating the skin from the meat. Rub the oil mixture all over the chicken, spooning some under the skin covering the breast. Tie the legs together with butchers twine and tuck the wings behind the back.
Step 2: Toss the sliced potatoes with the remaining tablespoon of oil, 3/4 t. salt and 1/2 t. pepper. Arrange the potatoes in a 12-inch nonstick, oven-safe skillet (see note above).
Step 3: Place the skillet over medium heat and cook, without flipping or stirring until the potatoes are golden-brown on the bottoms, about 7-9 minutes. Then place the prepared chicken, breast-side up on top of the potatoes.
Step 4: Transfer the skillet to the oven and roast until a thermometer inserted in several places has reached 165 degrees (about 1-1 1/4 hours).
Step 5: Transfer the chicken to a serving platter and cover loosely with foil, let rest for 20 minutes.
Step 6: Meanwhile, cover the skillet with the potatoes, add back to the oven and roast for an additional 20 minutes, until the potatoes are very tender. Carve the chicken and serve along side the roasted potatoes (and prepared to be wowed!)."
$ws.Range("A4").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Rows(15).AutoFit()

$excel.CutCopyMode = $false
